$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 13136
$ws.Range("E2").Value = 412
$ws.Range("F2").Value = 412
$ws.Range("G2").Value = 233
$ws.Range("H2").Value = 439
$ws.Range("I2").Value = 440
$ws.Range("J2").Value = -1
$ws.Range("K2").Value = 9275
$ws.Range("L2").Value = 5499
$ws.Range("M2").Value = 3776
$ws.Range("N2").Value = 3759
$ws.Range("O2").Value = 17
$ws.Range("P2").Value = 175
$ws.Range("Q2").Value = 700
$ws.Range("R2").Value = -922
$ws.Range("S2").Value = 901
$ws.Range("T2").Value = 874
$ws.Range("U2").Value = -174
$ws.Range("V2").Value = 3220
$ws.Range("W2").Value = 3.13
$ws.Range("X2").Value = 3.34
$ws.Range("Y2").Value = 12.41
$ws.Range("Z2").Value = 5.14
$ws.Range("AA2").Value = 145.62
$ws.Range("AB2").Value = 2130.88
$ws.Range("AC2").Value = 1259
$ws.Range("AD2").Value = 5.78
$ws.Range("AE2").Value = 10825
$ws.Range("AF2").Value = 0.67
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 1.37
$ws.Range("AI2").Value = 7.9
$ws.Range("AJ2").Value = 34920410

# Row 3
$ws.Range("D3").Value = 12517
$ws.Range("E3").Value = 192
$ws.Range("F3").Value = 192
$ws.Range("G3").Value = -138
$ws.Range("H3").Value = -246
$ws.Range("I3").Value = -232
$ws.Range("J3").Value = -14
$ws.Range("K3").Value = 9005
$ws.Range("L3").Value = 5371
$ws.Range("M3").Value = 3634
$ws.Range("N3").Value = 3628
$ws.Range("O3").Value = 6
$ws.Range("P3").Value = 175
$ws.Range("Q3").Value = 291
$ws.Range("R3").Value = -831
$ws.Range("S3").Value = -177
$ws.Range("T3").Value = 902
$ws.Range("U3").Value = -611
$ws.Range("V3").Value = 3123
$ws.Range("W3").Value = 1.53
$ws.Range("X3").Value = -1.96
$ws.Range("Y3").Value = -6.28
$ws.Range("Z3").Value = -2.69
$ws.Range("AA3").Value = 147.83
$ws.Range("AB3").Value = 1979.82
$ws.Range("AC3").Value = -664
$ws.Range("AD3").Value = -9.52
$ws.Range("AE3").Value = 10447
$ws.Range("AF3").Value = 0.6
$ws.Range("AG3").Value = 100
$ws.Range("AH3").Value = 1.58
$ws.Range("AI3").Value = -14.98
$ws.Range("AJ3").Value = 34920410

# Row 4
$ws.Range("D4").Value = 12497
$ws.Range("E4").Value = 462
$ws.Range("F4").Value = 462
$ws.Range("G4").Value = 590
$ws.Range("H4").Value = 482
$ws.Range("I4").Value = 478
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = 9641
$ws.Range("L4").Value = 5703
$ws.Range("M4").Value = 3938
$ws.Range("N4").Value = 3930
$ws.Range("O4").Value = 8
$ws.Range("P4").Value = 175
$ws.Range("Q4").Value = 576
$ws.Range("R4").Value = -354
$ws.Range("S4").Value = 32
$ws.Range("T4").Value = 469
$ws.Range("U4").Value = 107
$ws.Range("V4").Value = 3260
$ws.Range("W4").Value = 3.69
$ws.Range("X4").Value = 3.86
$ws.Range("Y4").Value = 12.64
$ws.Range("Z4").Value = 5.17
$ws.Range("AA4").Value = 144.8
$ws.Range("AB4").Value = 2236.26
$ws.Range("AC4").Value = 1368
$ws.Range("AD4").Value = 4.72
$ws.Range("AE4").Value = 11318
$ws.Range("AF4").Value = 0.57
$ws.Range("AG4").Value = 100
$ws.Range("AH4").Value = 1.55
$ws.Range("AI4").Value = 7.27
$ws.Range("AJ4").Value = 34920410

# Row 5
$ws.Range("D5").Value = 10669
$ws.Range("E5").Value = -228
$ws.Range("F5").Value = -228
$ws.Range("G5").Value = -478
$ws.Range("H5").Value = -309
$ws.Range("I5").Value = -303
$ws.Range("J5").Value = -6
$ws.Range("K5").Value = 8464
$ws.Range("L5").Value = 4937
$ws.Range("M5").Value = 3527
$ws.Range("N5").Value = 3527
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 175
$ws.Range("Q5").Value = -201
$ws.Range("R5").Value = -352
$ws.Range("S5").Value = 112
$ws.Range("T5").Value = 576
$ws.Range("U5").Value = -777
$ws.Range("V5").Value = 3187
$ws.Range("W5").Value = -2.14
$ws.Range("X5").Value = -2.9
$ws.Range("Y5").Value = -8.130000000000001
$ws.Range("Z5").Value = -3.42
$ws.Range("AA5").Value = 139.95
$ws.Range("AB5").Value = 2039.07
$ws.Range("AC5").Value = -868
$ws.Range("AD5").Value = -4.95
$ws.Range("AE5").Value = 10157
$ws.Range("AF5").Value = 0.42
$ws.Range("AG5").Value = 50
$ws.Range("AH5").Value = 1.16
$ws.Range("AI5").Value = -5.43
$ws.Range("AJ5").Value = 34920410

# Row 6
$ws.Range("D6").Value = 10758
$ws.Range("E6").Value = -211
$ws.Range("F6").Value = -211
$ws.Range("G6").Value = -475
$ws.Range("H6").Value = -571
$ws.Range("I6").Value = -571
$ws.Range("K6").Value = 8282
$ws.Range("L6").Value = 5350
$ws.Range("M6").Value = 2932
$ws.Range("N6").Value = 2932
$ws.Range("P6").Value = 175
$ws.Range("Q6").Value = 72
$ws.Range("R6").Value = -524
$ws.Range("S6").Value = 577
$ws.Range("T6").Value = 661
$ws.Range("U6").Value = -589
$ws.Range("V6").Value = 3808
$ws.Range("W6").Value = -1.96
$ws.Range("X6").Value = -5.31
$ws.Range("Y6").Value = -17.68
$ws.Range("Z6").Value = -6.82
$ws.Range("AA6").Value = 182.5
$ws.Range("AB6").Value = 1678.68
$ws.Range("AC6").Value = -1635
$ws.Range("AD6").Value = -1.44
$ws.Range("AE6").Value = 8904
$ws.Range("AF6").Value = 0.27
$ws.Range("AG6").Value = 25
$ws.Range("AH6").Value = 1.06
$ws.Range("AI6").Value = -1.44
$ws.Range("AJ6").Value = 34920410

# Row 7
$ws.Range("D7").Value = 11830
$ws.Range("E7").Value = 145
$ws.Range("G7").Value = 150
$ws.Range("I7").Value = 155
$ws.Range("W7").Value = 1.23
$ws.Range("AC7").Value = 444
$ws.Range("AD7").Value = 7.46
$ws.Range("AA7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()

# Row 8
$ws.Range("D8").Value = 12945
$ws.Range("E8").Value = 360
$ws.Range("G8").Value = 380
$ws.Range("I8").Value = 250
$ws.Range("W8").Value = 2.78
$ws.Range("AC8").Value = 716
$ws.Range("AD8").Value = 4.62
$ws.Range("AA8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()

# Row 9
$ws.Range("D9").Value = 13700
$ws.Range("E9").Value = 550
$ws.Range("G9").Value = 500
$ws.Range("I9").Value = 420
$ws.Range("W9").Value = 4.01
$ws.Range("AC9").Value = 1203
$ws.Range("AD9").Value = 2.75
$ws.Range("AA9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
